$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row with P-column SUMPRODUCT (old row 82) before shifting other rows
$ws.Rows("82:82").Delete()

# Remove the "theme testing" block of rows (old rows 56-76)
$ws.Rows("56:76").Delete()

# Remove now-unused column P
$ws.Columns("P:P").Delete()

# Fix the AutoFilter range to stop at row 56 (the blank separator row)
$ws.AutoFilterMode = $false
$ws.Rows("57:200").Insert()
$ws.Range("A1:K56").AutoFilter()
$ws.Rows("57:200").Delete()

# Fix the defined name ranges for _FilterDatabase
foreach ($dn in $wb.Names) {
    $dn.RefersTo = "=Sheet1!`$A`$1:`$K`$56"
}

Write-Host "Final UsedRange:" $ws.UsedRange.Address()
Write-Host "Final AutoFilter:" $ws.AutoFilter.Range.Address()
